# =====================================================================
# Applies the commit "Added removal of elements related to deleted
# things" to the "Do poprawienia:" bullet list (numId=18):
#
#   1. "Route Details - ..."                         (unchanged)
#   2. "Destinations Create - ..." [+ _GoBack bookmark]
#        -> text replaced by "Sprawienie, by ladowanie elementow do
#           tabelek odpowiednio dzialalo" + extra run " - Angular
#           Material?"; bookmark removed from here
#   3. "Warehouse - Delete - ..."
#        -> text replaced by "Destination - Create - generuje sie
#           nieprawidlowa kolejnosc czasami, a czasami prawidlowa -
#           nie mam pojecia dlaczego"; _GoBack bookmark placed here
#   4. "Route - Delete - ..."                        -> paragraph removed
#   5. "Deletions - ..."                              -> paragraph removed
#   6. "Uporzadkowanie kodu"                          -> paragraph removed
# =====================================================================

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the three paragraphs that are dropped completely.
# They are paragraphs 5, 6 and 7 in the original document. Delete from
# the bottom up so paragraph indices above stay valid while we work.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(7).Range.Delete()   # "Uporzadkowanie kodu"
$d.Paragraphs.Item(6).Range.Delete()   # "Deletions - ..."
$d.Paragraphs.Item(5).Range.Delete()   # "Route - Delete - ..."

# ---------------------------------------------------------------------
# Step 2: paragraph 4 ("Warehouse - Delete - ...") becomes the new
# "Destination - Create - ..." bullet, and gains the _GoBack bookmark
# that used to sit on paragraph 3.
# ---------------------------------------------------------------------

# Remove the old bookmark first (it currently lives at the end of
# paragraph 3's text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$p4 = $d.Paragraphs.Item(4)
$p4Text = $p4.Range
$p4Text.End = $p4Text.End - 1   # exclude the paragraph mark
$p4Text.Text = "Destination – Create – generuje się nieprawidłowa kolejność czasami, a czasami prawidłowa – nie mam pojęcia dlaczego"

# Re-fetch the (now shorter/longer) paragraph and append a one
# character placeholder. A bookmark cannot be reliably inserted with a
# zero-length range landing exactly on the last character of a
# paragraph, so we temporarily extend the paragraph by one character,
# anchor the bookmark just before that extra character, and then
# delete it again - the bookmark (being to the left of the deleted
# text) stays correctly anchored at the true end of the text.
$p4 = $d.Paragraphs.Item(4)
$p4Text = $p4.Range
$p4Text.End = $p4Text.End - 1
$p4Text.InsertAfter("X")

$p4 = $d.Paragraphs.Item(4)
$paraEnd = $p4.Range.End           # includes the paragraph mark
$bmPos = $paraEnd - 2              # right before the placeholder "X"
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($paraEnd - 2, $paraEnd - 1)
$placeholder.Delete()

# ---------------------------------------------------------------------
# Step 3: paragraph 3 ("Destinations Create - ...") becomes the new
# "Sprawienie, by ladowanie ..." bullet, written as two runs (the
# second one starting with a space, " - Angular Material?").
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3Text = $p3.Range
$p3Text.End = $p3Text.End - 1   # exclude the paragraph mark
$p3Text.Text = "Sprawienie, by ładowanie elementów do tabelek odpowiednio działało"

$p3 = $d.Paragraphs.Item(3)
$p3End = $p3.Range.Duplicate
$p3End.End = $p3End.End - 1
$p3End.Collapse(0)   # wdCollapseEnd
$p3End.InsertAfter(" – Angular Material?")
# Force the newly inserted text to live in its own run (rather than
# being merged back into the previous run) by toggling a character
# property on and back off.
$p3End.Bold = 1
$p3End.Bold = 0
